$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 210, shifting existing rows 210..311 down to 211..312
$ws.Rows.Item(210).Insert()

# Populate the newly inserted row 210 with the new record
$ws.Range("A210").Value = 11
$ws.Range("B210").Value = "Vega Monumental Concepción"
$ws.Range("C210").Value = "Bíobío"
$ws.Range("D210").Value = 45134
$ws.Range("D210").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E210").Value = 8
$ws.Range("F210").Value = 100112003
$ws.Range("G210").Value = "Ajo"
$ws.Range("H210").Value = "Chino"
$ws.Range("I210").Value = "Primera"
$ws.Range("J210").Value = 200
$ws.Range("K210").Value = 18000
$ws.Range("L210").Value = 19000
$ws.Range("M210").Value = 18500
$ws.Range("N210").Value = "$/caja 10 kilos"
$ws.Range("O210").Value = "China"
$ws.Range("P210").Value = 1850
$ws.Range("Q210").Value = 10
$ws.Range("R210").Value = "Hortaliza"
